# Apply cryptos list price/volume updates per commit "Updated cryptos list on Mon Apr 24 04:29:05 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.025.76"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "1.885.10"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.98%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.48"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4743"
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3947"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.00"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08011"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.019"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.87"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.901.00"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.054"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.201"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.014"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.51"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06734"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.10"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "28.010.53"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.509"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.347"
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").Value = "2.134.94"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.32"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.97"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.105"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.507"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.67"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9788"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09564"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.632"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.348"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.363"
$ws.Range("E36").Value = "  -5.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02253"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06070"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.206"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.203"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.010"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5974"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1889"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.33"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.265"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5660"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.934"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.348"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06813"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.32"
$ws.Range("E51").Value = "  -1.79%  "
